$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column D: "SPV" -> "Folio No"
$ws.Range("D1").Value = "Folio No"

# Replace SPV codes (AAA..EEE) with numeric folio numbers 1..5
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# Remove the Phone column (E) entirely
$ws.Range("E1:E6").ClearContents()

# Update the selected/active cell to D7 (matches saved selection state)
$ws.Range("D7").Select()
